$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Program" column (E) values for Term 1 rows so that labs (BK)
# are scheduled before classes (PM/FS/etc). This rewrites the existing
# shared-string text values in column E for rows 4, 6, 7, 9-14, 18-22.
$ws.Range("E4").Value = "BK"
$ws.Range("E6").Value = "PM"
$ws.Range("E7").Value = "PM"
$ws.Range("E9").Value = "BK"
$ws.Range("E10").Value = "PM"
$ws.Range("E11").Value = "PM"
$ws.Range("E12").Value = "PM"
$ws.Range("E13").Value = "PM"
$ws.Range("E14").Value = "BK"
$ws.Range("E18").Value = "BK"
$ws.Range("E19").Value = "BK"
$ws.Range("E20").Value = "BK"
$ws.Range("E21").Value = "BK"
$ws.Range("E22").Value = "BK"

# Update the active selection on the sheet to match the new selection.
$ws.Range("E10:E11").Select()
